# Auto-generated edit script: updates crypto price/volume table cells
# (and swaps two pairs of rows: 37/38 and 49/50) to match the
# "Updated cryptos list" GitHub Actions commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.306.87"
$ws.Range("E2").Value = "  -3.23%  "
$ws.Range("D3").Value = "2.468.34"
$ws.Range("E3").Value = "  -2.42%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'311.01"
$ws.Range("E5").Value = "  +1.07%  "
$ws.Range("D6").Value = "'94.36"
$ws.Range("E6").Value = "  -6.76%  "
$ws.Range("D7").Value = "'0.551"
$ws.Range("E7").Value = "  -2.86%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").Value = "'0.500"
$ws.Range("E9").Value = "  -5.02%  "
$ws.Range("D10").Value = "'33.41"
$ws.Range("E10").Value = "  -6.86%  "
$ws.Range("D11").Value = "'0.0778"
$ws.Range("E11").Value = "  -3.54%  "
$ws.Range("E12").Value = "  -0.54%  "
$ws.Range("E13").Value = "  -4.87%  "
$ws.Range("D14").Value = "2.844.13"
$ws.Range("E14").Value = "  -2.45%  "
$ws.Range("D15").Value = "2.477.54"
$ws.Range("E15").Value = "  -6.14%  "
$ws.Range("D16").Value = "'14.64"
$ws.Range("E16").Value = "  -6.72%  "
$ws.Range("D17").Value = "'0.783"
$ws.Range("E17").Value = "  -4.03%  "
$ws.Range("D18").Value = "41.239.70"
$ws.Range("E19").Value = "  -7.26%  "
$ws.Range("D20").Value = "0.0₃0917"
$ws.Range("E20").Value = "  -3.85%  "
$ws.Range("D21").Value = "'11.34"
$ws.Range("E21").Value = "  -6.95%  "
$ws.Range("D22").Value = "'68.25"
$ws.Range("E22").Value = "  -1.74%  "
$ws.Range("D23").Value = "'236.88"
$ws.Range("E23").Value = "  -2.96%  "
$ws.Range("D24").Value = "'2.77"
$ws.Range("E24").Value = "  -4.39%  "
$ws.Range("D25").Value = "'1.92"
$ws.Range("E25").Value = "  -5.93%  "
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("D27").Value = "'24.22"
$ws.Range("E27").Value = "  -6.21%  "
$ws.Range("D28").Value = "'2.23"
$ws.Range("E28").Value = "  -4.48%  "
$ws.Range("D29").Value = "'9.64"
$ws.Range("E29").Value = "  -5.17%  "
$ws.Range("D30").Value = "'35.80"
$ws.Range("E30").Value = "  -8.22%  "
$ws.Range("D31").Value = "'152.05"
$ws.Range("E31").Value = "  -3.65%  "
$ws.Range("D32").Value = "'5.55"
$ws.Range("E32").Value = "  -4.58%  "
$ws.Range("D33").Value = "'2.63"
$ws.Range("E33").Value = "  -5.99%  "
$ws.Range("D34").Value = "'2.59"
$ws.Range("E34").Value = "  -1.89%  "
$ws.Range("D35").Value = "'0.0743"
$ws.Range("E35").Value = "  -5.55%  "
$ws.Range("D36").Value = "'3.00"
$ws.Range("E36").Value = "  -6.65%  "
$ws.Range("B37").Value = "Celestia"
$ws.Range("C37").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D37").Value = "'17.24"
$ws.Range("E37").Value = "  -6.09%  "
$ws.Range("B38").Value = "ARBITRUM"
$ws.Range("C38").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D38").Value = "'1.89"
$ws.Range("E38").Value = "  -7.04%  "
$ws.Range("D39").Value = "'0.103"
$ws.Range("E39").Value = "  -7.60%  "
$ws.Range("E40").Value = "  -3.79%  "
$ws.Range("D41").Value = "'4.27"
$ws.Range("E41").Value = "  +1.54%  "
$ws.Range("D42").Value = "'20.54"
$ws.Range("E42").Value = "  -7.19%  "
$ws.Range("E43").Value = "  +0.15%  "
$ws.Range("D44").Value = "1.982.38"
$ws.Range("E44").Value = "  -0.89%  "
$ws.Range("D45").Value = "'0.0284"
$ws.Range("E45").Value = "  -5.70%  "
$ws.Range("D46").Value = "'3.05"
$ws.Range("E46").Value = "  -6.93%  "
$ws.Range("D47").Value = "'8.68"
$ws.Range("E47").Value = "  -2.83%  "
$ws.Range("D48").Value = "'70.00"
$ws.Range("E48").Value = "  -3.79%  "
$ws.Range("B49").Value = "Aave"
$ws.Range("C49").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D49").Value = "'96.54"
$ws.Range("E49").Value = "  -4.69%  "
$ws.Range("B50").Value = "BitcoinSV"
$ws.Range("C50").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D50").Value = "'75.00"
$ws.Range("E50").Value = "  -6.52%  "
$ws.Range("D51").Value = "'0.178"
$ws.Range("E51").Value = "  -6.79%  "
